$d = $word.ActiveDocument
$sec = $d.Sections.First

# --- Rename the two header logo shapes (swap image1.png / image2.png) ---
# Header 1 (primary/default header) currently holds the shape named "image2.png"
# and should become "image1.png".
$hdr1 = $sec.Headers.Item(1)
if ($hdr1.Exists) {
    for ($i = 1; $i -le $hdr1.Shapes.Count; $i++) {
        $shp = $hdr1.Shapes.Item($i)
        if ($shp.Name -eq "image2.png") {
            $shp.Name = "image1.png"
        }
    }
}

# Header 2 (first-page header) currently holds the shape named "image1.png"
# and should become "image2.png".
$hdr2 = $sec.Headers.Item(2)
if ($hdr2.Exists) {
    for ($i = 1; $i -le $hdr2.Shapes.Count; $i++) {
        $shp = $hdr2.Shapes.Item($i)
        if ($shp.Name -eq "image1.png") {
            $shp.Name = "image2.png"
        }
    }

    # --- Rename the mission title text "Sapling-1" -> "Sapling" ---
    $rng = $hdr2.Range.Duplicate
    $found = $rng.Find.Execute("Sapling-1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = "Sapling"
    }
}
